$wb = $excel.ActiveWorkbook

# ALC!row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6543.381
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 6543.381
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 19630.143
$ws.Range("N17").Value = -19966.143

# ALC!row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 252125
$ws.Range("I64").Value = 1000000
$ws.Range("J64").Value = 2833.3333
$ws.Range("K64").Value = 1000000
$ws.Range("L64").Value = 2833.3333
$ws.Range("M64").Value = -999752
$ws.Range("N64").Value = -3329.3333

# ALC!row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 252125
$ws.Range("I67").Value = 1000000
$ws.Range("J67").Value = 2833.3333
$ws.Range("K67").Value = 1000000
$ws.Range("L67").Value = 2833.3333
$ws.Range("M67").Value = -999142
$ws.Range("N67").Value = -4549.3333

# ALC!row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# ALC!row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# ALC!row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = 0

# ALC!row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = 0

# ALC!row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# ALC!row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 45900
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 45900
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 45900
$ws.Range("N130").Value = -55940

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19997.945
$ws.Range("I32").Value = 21009.488
$ws.Range("J32").Value = 9376.75
$ws.Range("K32").Value = 21009.488
$ws.Range("L32").Value = 9376.75
$ws.Range("M32").Value = -20722.488
$ws.Range("N32").Value = -9950.75

# ARM!row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 49991
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 49991
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 49991
$ws.Range("N80").Value = -51987

# ARM!row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 49991
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 49991
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 149973
$ws.Range("N83").Value = -159957

# ARM!row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 43371
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 43371
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 43371
$ws.Range("N109").Value = -46145

# ARM!row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 51707
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 51707
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 51707
$ws.Range("N131").Value = -61787

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7814335
$ws.Range("I132").Value = 11629121
$ws.Range("J132").Value = 3107.3809
$ws.Range("K132").Value = 34887363
$ws.Range("L132").Value = 9322.1427
$ws.Range("M132").Value = -34884833
$ws.Range("N132").Value = -14382.1427

# BSM!row 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()

# BSM!row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3845.6428
$ws.Range("I80").Value = 16775
$ws.Range("J80").Value = 319.45456
$ws.Range("K80").Value = 16775
$ws.Range("L80").Value = 319.45456
$ws.Range("M80").Value = -15777
$ws.Range("N80").Value = -2315.45456

# BSM!row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 3845.6428
$ws.Range("I83").Value = 16775
$ws.Range("J83").Value = 319.45456
$ws.Range("K83").Value = 83875
$ws.Range("L83").Value = 1597.2728
$ws.Range("M83").Value = -78883
$ws.Range("N83").Value = -11581.2728

# BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2678.647
$ws.Range("I105").Value = 2610.0667
$ws.Range("J105").Value = 3193
$ws.Range("K105").Value = 2610.0667
$ws.Range("L105").Value = 3193
$ws.Range("M105").Value = -863.0666999999999
$ws.Range("N105").Value = -6687

# BSM!row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 48477.855
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48477.855
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48477.855
$ws.Range("N130").Value = -58517.855

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4630.403
$ws.Range("I31").Value = 1898.9642
$ws.Range("J31").Value = 6591.436
$ws.Range("K31").Value = 1898.9642
$ws.Range("L31").Value = 6591.436
$ws.Range("M31").Value = -1603.9642
$ws.Range("N31").Value = -7181.436

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4630.403
$ws.Range("I34").Value = 1898.9642
$ws.Range("J34").Value = 6591.436
$ws.Range("K34").Value = 1898.9642
$ws.Range("L34").Value = 6591.436
$ws.Range("M34").Value = -1696.9642
$ws.Range("N34").Value = -6995.436

# CUL!row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 777.5
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 870
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 2610
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -3602

# CUL!row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5940.357
$ws.Range("I139").Value = 7573.3887
$ws.Range("J139").Value = 3000.9
$ws.Range("K139").Value = 22720.1661
$ws.Range("L139").Value = 9002.700000000001
$ws.Range("M139").Value = -17580.1661
$ws.Range("N139").Value = -19282.7

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 459544.53
$ws.Range("I80").Value = 504699
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 504699
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -503701
$ws.Range("N80").Value = -9996

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 459544.53
$ws.Range("I83").Value = 504699
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 2523495
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -2518503
$ws.Range("N83").Value = -49984

# GSM!row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 44728
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 44728
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 44728
$ws.Range("N130").Value = -54768

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2198.5
$ws.Range("I122").Value = 2042.1428
$ws.Range("J122").Value = 2563.3333
$ws.Range("K122").Value = 6126.428400000001
$ws.Range("L122").Value = 7689.999899999999
$ws.Range("M122").Value = -3676.428400000001
$ws.Range("N122").Value = -12589.9999

# WVR!row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2043.2307
$ws.Range("I81").Value = 2008.3334
$ws.Range("J81").Value = 2073.1428
$ws.Range("K81").Value = 4016.6668
$ws.Range("L81").Value = 4146.2856
$ws.Range("M81").Value = -2955.6668
$ws.Range("N81").Value = -6268.2856

# WVR!row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2043.2307
$ws.Range("I84").Value = 2008.3334
$ws.Range("J84").Value = 2073.1428
$ws.Range("K84").Value = 20083.334
$ws.Range("L84").Value = 20731.428
$ws.Range("M84").Value = -14779.334
$ws.Range("N84").Value = -31339.428

# WVR!row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1251
$ws.Range("I113").Value = 1251
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3753
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1583

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20588.982
$ws.Range("I136").Value = 65527.35
$ws.Range("J136").Value = 2399.6428
$ws.Range("K136").Value = 196582.05
$ws.Range("L136").Value = 7198.928400000001
$ws.Range("M136").Value = -194032.05
$ws.Range("N136").Value = -12298.9284
